$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D1 was is_locked_lbl -> becomes order_by
$ws.Range("D1").Value = '<%=comment.order_by%>'

# E1 was is_enabled_lbl -> becomes rem
$ws.Range("E1").Value = '<%=comment.rem%>'

# F1 was order_by -> becomes the new tenant_id_lbl validation string
$ws.Range("F1").Value = '<%=comment.tenant_id_lbl%><%selectList.tenant_id = data.findAllTenant.map((item) => item.lbl)%><%_dataValidation_({ sqref: `${ _col }2:${ _col }${ _lastRow }`, formula1: `"${ selectList.tenant_id.join(",") }"` })%>'

# G1 was rem -> cell removed entirely
$ws.Range("G1").ClearContents()
